$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.530.30'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.638.92'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''212.97'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').Value = '  +4.58%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''22.92'
$ws.Range('E8').Value = '  -5.53%  '
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').Value = '''0.0889'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('D12').Value = '1.872.58'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = '1.636.89'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').Value = '''4.02'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '''0.564'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '''64.16'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').Value = '27.547.17'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '''229.19'
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').Value = '''7.70'
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D22').Value = '''4.30'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('D23').Value = '''9.88'
$ws.Range('E23').Value = '  +5.95%  '
$ws.Range('E24').Value = '  -3.84%  '
$ws.Range('D25').Value = '''149.66'
$ws.Range('E25').Value = '  +2.13%  '
$ws.Range('D26').Value = '''6.96'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '''15.57'
$ws.Range('E29').Value = '  -2.67%  '
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').Value = '''0.0487'
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('E32').Value = '  -0.77%  '
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').Value = '1.423.81'
$ws.Range('E34').Value = '  -2.79%  '
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '''0.571'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '''0.877'
$ws.Range('E38').Value = '  -3.95%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = '''0.901'
$ws.Range('E40').Value = '  +14.97%  '
$ws.Range('E41').Value = '  -0.80%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '''2.26'
$ws.Range('D46').Value = '''64.89'
$ws.Range('E46').Value = '  -1.39%  '
$ws.Range('D47').Value = '1.781.63'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('D49').Value = '''86.16'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '''0.0986'
$ws.Range('E51').Value = '  -2.57%  '
